$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AX (col 50) -- shifts AX..BG right to AY..BH
$ws.Columns("AX").Insert()

# Populate the newly inserted column AX with the "Y23-ticket" data
$ws.Range("AX1").Value = "Y23-ticket"
$ws.Range("AX2").Value = 5
$ws.Range("AX3").Value = "2;3"

# New width for column AR (44)
$ws.Columns("AR").ColumnWidth = 24

# Update view: scroll position and selection
$excel.ActiveWindow.ScrollColumn = 28
$ws.Range("AS1:AS1048576").Select()
